$d = $word.ActiveDocument

# Remove the trailing paragraphs that are no longer needed (from the end,
# so indices of the earlier paragraphs stay valid while we work).
$d.Paragraphs.Item(5).Range.Delete()
$d.Paragraphs.Item(4).Range.Delete()
$d.Paragraphs.Item(3).Range.Delete()

# Build the single paragraph that should remain: a field (begin / instrText
# " " + "m:" + "endfor " / end) followed by a bold red run reporting the
# new error message.
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:instrText xml:space="preserve"> </w:instrText></w:r><w:r><w:instrText>m:</w:instrText></w:r><w:r><w:instrText xml:space="preserve">endfor </w:instrText></w:r><w:r><w:fldChar w:fldCharType="end"/></w:r><w:r><w:rPr><w:b w:val="true"/><w:color w:val="FF0000"/></w:rPr><w:t>Invalid if statement: Unexpected tag m:endfor at this location</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# Insert it as a brand-new paragraph at the very end of the document (this
# becomes the document's new final paragraph mark), then drop the original
# two leading paragraphs that used to hold the old content.
$endRange = $d.Range($d.Content.End, $d.Content.End)
$endRange.InsertXML($xml)

$d.Paragraphs.Item(2).Range.Delete()
$d.Paragraphs.Item(1).Range.Delete()
